# Update "想去人数" (number of people interested) figures on the
# "展览" and "全部类型" sheets:
#   F2: 324 -> 325
#   F4: 58  -> 60
#   F5: 282 -> 283

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 325
    $ws.Range("F4").Value = 60
    $ws.Range("F5").Value = 283
}
